# Weekly price update: insert a new "Haba" record for
# "Terminal La Palmera de La Serena" at row 29, pushing the existing
# rows 29-35 down to 30-36 (their values are unchanged, only their row
# position shifts).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 29; this shifts rows 29:35 -> 30:36.
$ws.Rows.Item(29).Insert()

# Populate the new row 29 with this week's data.
$ws.Cells.Item(29, 1).Value = 8
$ws.Cells.Item(29, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(29, 3).Value = "Coquimbo"
$ws.Cells.Item(29, 4).Value = 44855
$ws.Cells.Item(29, 5).Value = 4
$ws.Cells.Item(29, 6).Value = 100112026
$ws.Cells.Item(29, 7).Value = "Haba"
$ws.Cells.Item(29, 8).Value = "Sin especificar"
$ws.Cells.Item(29, 9).Value = "Primera"
$ws.Cells.Item(29, 10).Value = 540
$ws.Cells.Item(29, 11).Value = 7000
$ws.Cells.Item(29, 12).Value = 8000
$ws.Cells.Item(29, 13).Value = 7500
$ws.Cells.Item(29, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(29, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(29, 16).Value = 300
$ws.Cells.Item(29, 17).Value = 25
$ws.Cells.Item(29, 18).Value = "Hortaliza"
